$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the cell fill/styles that changed (black fill = "done" style used
#    elsewhere in the sheet as s=4, white fill = s=6, bordered/no-fill = s=1).
#    We reproduce the exact look by copying formats (not numeric style ids)
#    from stable donor cells that already carry the desired style:
#      B3  -> plain black-filled ("not yet" / s=4) cell
#      B10 -> white-filled text cell (s=6)
#      C5  -> bordered, no-fill text cell (s=1)
# ---------------------------------------------------------------------------

$blackFillDonor = $ws.Range("B3")     # style s=4
$whiteFillDonor = $ws.Range("B10")    # style s=6

# H12 goes from white(6) to black(4)
$blackFillDonor.Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null

# H13, F14, H14 go from black(4, empty) to white(6, with text)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$whiteFillDonor.Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$whiteFillDonor.Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null

# B31 goes from black(4) to white(6)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null

# B32 goes from black(4, empty) to white(6, with text)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null

# B34 goes from bordered/no-fill(1) to black(4)
$blackFillDonor.Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4122) | Out-Null

# B36, B37 go from black(4) to white(6)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null

# B38 goes from black(4, empty) to white(6, with text)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B38").PasteSpecial(-4122) | Out-Null

# B40 goes from bordered/no-fill(1) to black(4)
$blackFillDonor.Copy() | Out-Null
$ws.Range("B40").PasteSpecial(-4122) | Out-Null

# B42 goes from black(4, empty) to white(6, with text)
$whiteFillDonor.Copy() | Out-Null
$ws.Range("B42").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Update the cell text contents.
# ---------------------------------------------------------------------------

$ws.Range("G12").Value2 = "estudiar cálculo integral"
$ws.Range("G13").Value2 = "estudiar cálculo integral"
$ws.Range("H13").Value2 = "estudiar física mecánica"
$ws.Range("F14").Value2 = "estudiar/prácticas/tareas informática"
$ws.Range("H14").Value2 = "estudiar física mecánica"

$ws.Range("B32").Value2 = "estudiar física mecánica"

$ws.Range("C33").Value2 = "estudiar algebra lineal"

$ws.Range("B34").Value2 = "Adelantar prácticas laboratorio informatica"
$ws.Range("C34").Value2 = "estudiar algebra lineal"

$ws.Range("B35").Value2 = "Adelantar prácticas laboratorio informatica"
$ws.Range("B36").Value2 = "Adelantar prácticas laboratorio informatica"
$ws.Range("B37").Value2 = "Adelantar prácticas laboratorio informatica"
$ws.Range("B38").Value2 = "Adelantar prácticas laboratorio informatica"

$ws.Range("B40").Value2 = "Estudiar algebra líneal"
$ws.Range("B41").Value2 = "Estudiar algebra líneal"
$ws.Range("B42").Value2 = "Estudiar algebra líneal"

# ---------------------------------------------------------------------------
# 3) Update the view: scroll so row 17 is at the top, and select C36
#    (best effort -- selection/activeCell persists; scroll position is
#    restored as closely as the host allows).
# ---------------------------------------------------------------------------

$ws.Range("A17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C36").Select() | Out-Null

"Done"
